$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append new data row 88 with the Q-factor run results for
# sg_rr_100_025 2023-12-11 14-23-14
$ws.Range("A88").Value = "sg_rr_100_025 2023-12-11 14-23-14.csv"
$ws.Range("B88").Value = 0.01
$ws.Range("C88").Value = 1000
$ws.Range("D88").Value = 5001
$ws.Range("E88").Value = 1530
$ws.Range("F88").Value = 1570
$ws.Range("G88").Value = 0.001
$ws.Range("H88").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I88").Value = 1
$ws.Range("J88").Value = 0.98282051282051597
$ws.Range("K88").Value = 0.0055097596875867197
$ws.Range("L88").Value = "yes,although possible peak at end missed although hard to tell if it is a peak"
$ws.Range("M88").Value = 0.16184922620443801
$ws.Range("N88").Value = 0.0060009599444959304
$ws.Range("O88").Value = 9978.8525564459196
$ws.Range("P88").Value = 283.12224613964702
$ws.Range("Q88").Value = 99366917.390740097
$ws.Range("R88").Value = 8460116.0168236997
$ws.Range("S88").Value = 100
$ws.Range("T88").Value = 0.1

# Update the view: scroll position and active cell/selection moved down
# one row to reflect the newly added row.
$excel.ActiveWindow.ScrollRow = 72
$ws.Range("A89").Select() | Out-Null
